# Se completa la clase AsignacionNuevos | Se empieza con los tests de AsignacionNuevos
#
# The sheet originally has columns:
#   A cod_CNO | B recursosxcno | C n_programas
# It becomes:
#   A cod_CNO | B ipo_ponderado | C cuposxcno | D ipo | E n_programas
#   | F participacion_ipo | G recursosxcno

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new header cells and copy the existing header
#     formatting (bold font, thin border, centered alignment) onto them ---
$ws.Range("A1").Copy($ws.Range("C1"))
$ws.Range("A1").Copy($ws.Range("D1"))
$ws.Range("A1").Copy($ws.Range("E1"))
$ws.Range("A1").Copy($ws.Range("F1"))
$ws.Range("A1").Copy($ws.Range("G1"))

# --- Header row ---
$ws.Range("A1").Value = "cod_CNO"
$ws.Range("B1").Value = "ipo_ponderado"
$ws.Range("C1").Value = "cuposxcno"
$ws.Range("D1").Value = "ipo"
$ws.Range("E1").Value = "n_programas"
$ws.Range("F1").Value = "participacion_ipo"
$ws.Range("G1").Value = "recursosxcno"

# --- Data rows: cod_CNO, ipo_ponderado, cuposxcno, ipo, n_programas, participacion_ipo, recursosxcno ---
$data = @(
    @(1231, 96.55310409621252,  200, 0.9655310409621252, 2, 0.08149173855568445, 107569094.8935035),
    @(1345, 99.74631914374956,  160, 1.24682898929687,   2, 0.08418684243910454, 111126632.019618),
    @(2242, 40.66875408474579,   75, 0.5422500544632771, 1, 0.03432481540890753, 45308756.33975793),
    @(2281, 55.2338365446359,    83, 0.6654679101763361, 1, 0.04661788359116507, 61535606.34033789),
    @(3311, 47.73412360227283,  100, 0.4773412360227283, 1, 0.04028805450837633, 53180231.95105676),
    @(3315, 95.97906473146111,  160, 1.199738309143264,  2, 0.0810072438698753,  106929561.9082354),
    @(4311, 143.1970855311291,  300, 1.43197085531129,   3, 0.120859702702158,   159534807.5668485),
    @(5223, 15.08304680655133,   30, 0.5027682268850443, 1, 0.01273023501924728, 16803910.2254064),
    @(6234, 32.62370608853713,   60, 0.5437284348089521, 1, 0.02753471835183438, 36345828.22442139),
    @(6322, 158.1962925708634,  285, 1.110149421549919,  2, 0.1335191761604801,  176245312.5318337),
    @(6355, 194.070393492528,   350, 2.217947354200319,  4, 0.1637972586788359,  216212381.4560634),
    @(8325, 205.7350392827064,  410, 2.508963893691542,  5, 0.1736423307143312,  229207876.5429172)
)

$rowIndex = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex = $rowIndex + 1
}
